$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first data row (old row 2); everything shifts up, matching
# the new "traditional agent" offset: Trad_Prediction(row) = prior Real_Close.
$ws.Rows(2).Delete()

# Rebuild Trad_Prediction (C) and AI_Prediction (D) with the new simulation values.
$ws.Cells.Item(2, 3).Value = 209.9100036621094
$ws.Cells.Item(2, 4).Value = 207.0885515158442
$ws.Cells.Item(3, 3).Value = 207.8049926757812
$ws.Cells.Item(3, 4).Value = 207.3134535566308
$ws.Cells.Item(4, 3).Value = 207.8000030517578
$ws.Cells.Item(4, 4).Value = 198.0988443588153
$ws.Cells.Item(5, 3).Value = 208.9156036376953
$ws.Cells.Item(5, 4).Value = 213.1562751313305
$ws.Cells.Item(6, 3).Value = 209.7324981689453
$ws.Cells.Item(6, 4).Value = 209.9908352218939
$ws.Cells.Item(7, 3).Value = 209.8899993896484
$ws.Cells.Item(7, 4).Value = 209.5294136496121
$ws.Cells.Item(8, 3).Value = 211.1100006103516
$ws.Cells.Item(8, 4).Value = 208.0585043743041
$ws.Cells.Item(9, 3).Value = 212.0099945068359
$ws.Cells.Item(9, 4).Value = 209.5979425958947
$ws.Cells.Item(10, 3).Value = 212.9799957275391
$ws.Cells.Item(10, 4).Value = 209.9980258609183
$ws.Cells.Item(11, 3).Value = 212.0800018310547
$ws.Cells.Item(11, 4).Value = 214.8890764744741
$ws.Cells.Item(12, 3).Value = 212.6199951171875
$ws.Cells.Item(12, 4).Value = 207.067666099176
$ws.Cells.Item(13, 3).Value = 212.6000061035156
$ws.Cells.Item(13, 4).Value = 209.586717867964
$ws.Cells.Item(14, 3).Value = 212.5749969482422
$ws.Cells.Item(14, 4).Value = 214.3841571570057
$ws.Cells.Item(15, 3).Value = 212.4199981689453
$ws.Cells.Item(15, 4).Value = 214.0784842464163
$ws.Cells.Item(16, 3).Value = 211.7899932861328
$ws.Cells.Item(16, 4).Value = 207.416602760971
$ws.Cells.Item(17, 3).Value = 210.8249969482422
$ws.Cells.Item(17, 4).Value = 208.2199641111351
$ws.Cells.Item(18, 3).Value = 210.6349945068359
$ws.Cells.Item(18, 4).Value = 216.6839633550513
$ws.Cells.Item(19, 3).Value = 211.1649932861328
$ws.Cells.Item(19, 4).Value = 211.8546269106953
$ws.Cells.Item(20, 3).Value = 211.4049987792969
$ws.Cells.Item(20, 4).Value = 213.880629082019
$ws.Cells.Item(21, 3).Value = 211.1600036621094
$ws.Cells.Item(21, 4).Value = 206.4115685202648
$ws.Cells.Item(22, 3).Value = 211.0749969482422
$ws.Cells.Item(22, 4).Value = 212.7758078496427
$ws.Cells.Item(23, 3).Value = 208.0299987792969
$ws.Cells.Item(23, 4).Value = 211.7921301119123
$ws.Cells.Item(24, 3).Value = 208.7649993896484
$ws.Cells.Item(24, 4).Value = 207.6747298772758
$ws.Cells.Item(25, 3).Value = 208.8406982421875
$ws.Cells.Item(25, 4).Value = 211.6665997934498
$ws.Cells.Item(26, 3).Value = 209
$ws.Cells.Item(26, 4).Value = 208.2140065425759
$ws.Cells.Item(27, 3).Value = 209
$ws.Cells.Item(27, 4).Value = 212.8942053003424
$ws.Cells.Item(28, 3).Value = 209.1999053955078
$ws.Cells.Item(28, 4).Value = 214.5460316979326
$ws.Cells.Item(29, 3).Value = 208.5399932861328
$ws.Cells.Item(29, 4).Value = 213.3241945868067
$ws.Cells.Item(30, 3).Value = 210.1699981689453
$ws.Cells.Item(30, 4).Value = 206.1443109000217
$ws.Cells.Item(31, 3).Value = 211.4299926757812
$ws.Cells.Item(31, 4).Value = 212.4125054596048
$ws.Cells.Item(32, 3).Value = 210.6799926757812
$ws.Cells.Item(32, 4).Value = 211.4057365283323
$ws.Cells.Item(33, 3).Value = 210.6726989746094
$ws.Cells.Item(33, 4).Value = 216.4175418824649
$ws.Cells.Item(34, 3).Value = 211.375
$ws.Cells.Item(34, 4).Value = 208.5729180807568
$ws.Cells.Item(35, 3).Value = 210.1799926757812
$ws.Cells.Item(35, 4).Value = 211.8234553491991
